# Updated cryptos list on Thu Jan 11 12:13:45 UTC 2024 with GitHub Actions
# Refresh price/volume figures for the coin table on Sheet1 (rows 2-51).
# Note: numeric-looking price strings (e.g. "315.68") are written with a
# leading apostrophe so Excel keeps them as text instead of auto-converting
# them to numbers, matching the original inline-string cell contents.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.102.98"
$ws.Range("E2").Value = "  +3.31%  "
$ws.Range("D3").Value = "2.652.22"
$ws.Range("E3").Value = "  +10.14%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'315.68"
$ws.Range("E5").Value = "  +5.59%  "
$ws.Range("D6").Value = "'104.61"
$ws.Range("E6").Value = "  +7.54%  "
$ws.Range("D7").Value = "'0.613"
$ws.Range("E7").Value = "  +8.91%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.601"
$ws.Range("E9").Value = "  +17.22%  "
$ws.Range("D10").Value = "'39.67"
$ws.Range("E10").Value = "  +13.77%  "
$ws.Range("D11").Value = "'55.84"
$ws.Range("E11").Value = "  +3.70%  "
$ws.Range("D12").Value = "'0.0861"
$ws.Range("E12").Value = "  +8.75%  "
$ws.Range("D13").Value = "'8.47"
$ws.Range("E13").Value = "  +19.08%  "
$ws.Range("D14").Value = "3.058.09"
$ws.Range("E14").Value = "  +10.16%  "
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").Value = "2.677.02"
$ws.Range("E16").Value = "  +11.07%  "
$ws.Range("D17").Value = "'0.947"
$ws.Range("E17").Value = "  +12.81%  "
$ws.Range("D18").Value = "'15.40"
$ws.Range("E18").Value = "  +8.46%  "
$ws.Range("D19").Value = "47.682.87"
$ws.Range("E19").Value = "  +4.61%  "
$ws.Range("E20").Value = "  +10.01%  "
$ws.Range("D21").Value = "'13.51"
$ws.Range("E21").Value = "  +5.65%  "
$ws.Range("D22").Value = "'6.88"
$ws.Range("E22").Value = "  +11.17%  "
$ws.Range("D23").Value = "'72.97"
$ws.Range("E23").Value = "  +8.69%  "
$ws.Range("D24").Value = "'276.74"
$ws.Range("E24").Value = "  +14.62%  "
$ws.Range("D25").Value = "'3.13"
$ws.Range("E25").Value = "  +11.41%  "
$ws.Range("D26").Value = "'31.43"
$ws.Range("E26").Value = "  +47.43%  "
$ws.Range("E27").Value = "  +16.81%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "'4.07"
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("D30").Value = "'10.84"
$ws.Range("E30").Value = "  +11.25%  "
$ws.Range("D31").Value = "'41.81"
$ws.Range("E31").Value = "  +9.37%  "
$ws.Range("D32").Value = "'2.31"
$ws.Range("E32").Value = "  +4.04%  "
$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D33").Value = "'3.86"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'6.28"
$ws.Range("E34").Value = "  +13.96%  "
$ws.Range("D35").Value = "'2.33"
$ws.Range("E35").Value = "  +17.06%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.0855"
$ws.Range("E36").Value = "  +10.86%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'2.87"
$ws.Range("E37").Value = "  +5.70%  "
$ws.Range("D38").Value = "'153.77"
$ws.Range("E38").Value = "  +2.27%  "
$ws.Range("D39").Value = "'0.121"
$ws.Range("E39").Value = "  +7.06%  "
$ws.Range("D40").Value = "'0.125"
$ws.Range("E40").Value = "  +7.85%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "'24.46"
$ws.Range("E41").Value = "  +53.44%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'16.94"
$ws.Range("E42").Value = "  +12.02%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'4.36"
$ws.Range("E43").Value = "  +13.51%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "'3.74"
$ws.Range("E44").Value = "  +15.07%  "
$ws.Range("D45").Value = "'0.0335"
$ws.Range("E45").Value = "  +12.12%  "
$ws.Range("D46").Value = "2.086.96"
$ws.Range("E46").Value = "  +6.60%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "'98.99"
$ws.Range("E47").Value = "  +9.30%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "'0.999"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "'116.53"
$ws.Range("E49").Value = "  +14.27%  "
$ws.Range("D50").Value = "'1.89"
$ws.Range("E50").Value = "  +9.05%  "
$ws.Range("D51").Value = "'9.30"
$ws.Range("E51").Value = "  +6.27%  "
